# Update gh-pages output (杭州-漫展信息) workbook to the data scraped at 456a3b4.
#
# Sheet layout (1-indexed via Worksheets.Item):
#   1 = 展览       (exhibitions)
#   2 = 演出       (performances)
#   3 = 本地生活   (local life)
#   4 = 全部类型   (all types / combined feed)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 - "想去人数" (column F) count bumps
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 6).Value  = 272
$ws1.Cells.Item(3, 6).Value  = 3229
$ws1.Cells.Item(4, 6).Value  = 1997
$ws1.Cells.Item(5, 6).Value  = 267
$ws1.Cells.Item(6, 6).Value  = 103
$ws1.Cells.Item(7, 6).Value  = 3083
$ws1.Cells.Item(8, 6).Value  = 613
$ws1.Cells.Item(9, 6).Value  = 301
$ws1.Cells.Item(12, 6).Value = 152
$ws1.Cells.Item(15, 6).Value = 10137
$ws1.Cells.Item(19, 6).Value = 35
$ws1.Cells.Item(20, 6).Value = 8026
$ws1.Cells.Item(21, 6).Value = 12640
$ws1.Cells.Item(24, 6).Value = 21
$ws1.Cells.Item(26, 6).Value = 395
$ws1.Cells.Item(27, 6).Value = 593
$ws1.Cells.Item(29, 6).Value = 413
$ws1.Cells.Item(30, 6).Value = 2823
$ws1.Cells.Item(33, 6).Value = 7956
$ws1.Cells.Item(34, 6).Value = 1518
$ws1.Cells.Item(38, 6).Value = 4620
$ws1.Cells.Item(39, 6).Value = 1407
$ws1.Cells.Item(43, 6).Value = 638

# ---------------------------------------------------------------------------
# Sheet 2: 演出 - count bump on row 14, plus a new early-bird concert
# inserted as row 22 (pushing the old rows 22-24 down to 23-25).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(14, 6).Value = 3

# Insert a blank row at 22; the engine shifts the old rows 22-24 (and all of
# their formatting / text-vs-number typing) down to 23-25, including the
# "A" index column's existing numeric values (old A22=21 physically moves to
# A23, etc).
$ws2.Rows(22).Insert()

# New row 22 content: "加勒比海盗/权力的游戏" new year concert.
$ws2.Cells.Item(22, 1).Value = 21
# Column B holds plain "yyyy-mm-dd" text in this sheet; force text typing
# before assigning so it isn't auto-parsed into a date serial number.
$ws2.Cells.Item(22, 2).NumberFormat = "@"
$ws2.Cells.Item(22, 2).Value = "2025-01-01"
$ws2.Cells.Item(22, 3).Value = "杭州·【早鸟5折】2025超燃视听新年交响音乐会《加勒比海盗》《权力的游戏》"
$ws2.Cells.Item(22, 4).Value = "湖墅南路138号 杭州浙话艺术剧院"
$ws2.Cells.Item(22, 5).Value = "2025.01.01 20:00-01.01 21:30"
$ws2.Cells.Item(22, 6).Value = 0
$ws2.Cells.Item(22, 7).Value = 140
$ws2.Cells.Item(22, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93886"
$ws2.Cells.Item(22, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/qDzvMMvq1729667982109.jpeg"

# Tidy up formatting: drop the ad-hoc styles picked up from the insert /
# text-format dance, then restore the bold+bordered "index column" look on
# A22 (matching every other row's A cell) by copying the format from A21.
$ws2.Range("A22:I22").ClearFormats()
$ws2.Cells.Item(21, 1).Copy()
$ws2.Cells.Item(22, 1).PasteSpecial(-4122)

# The rows that got physically shifted down (old 22/23/24 -> new 23/24/25)
# kept their original numeric "A" index values from before the shift, so
# just renumber that column to match their new row position.
$ws2.Cells.Item(23, 1).Value = 22
$ws2.Cells.Item(24, 1).Value = 23
$ws2.Cells.Item(25, 1).Value = 24

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 - "想去人数" (column F) count bumps
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(4, 6).Value = 226
$ws3.Cells.Item(5, 6).Value = 19

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 - "想去人数" (column F) count bumps
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3, 6).Value  = 272
$ws4.Cells.Item(4, 6).Value  = 226
$ws4.Cells.Item(5, 6).Value  = 3229
$ws4.Cells.Item(7, 6).Value  = 1997
$ws4.Cells.Item(9, 6).Value  = 267
$ws4.Cells.Item(10, 6).Value = 19
$ws4.Cells.Item(11, 6).Value = 3083
$ws4.Cells.Item(13, 6).Value = 613
$ws4.Cells.Item(16, 6).Value = 152
$ws4.Cells.Item(19, 6).Value = 10137
$ws4.Cells.Item(22, 6).Value = 35
$ws4.Cells.Item(23, 6).Value = 8026
$ws4.Cells.Item(24, 6).Value = 12640
$ws4.Cells.Item(27, 6).Value = 21
$ws4.Cells.Item(30, 6).Value = 593
$ws4.Cells.Item(33, 6).Value = 2823
$ws4.Cells.Item(38, 6).Value = 7956
$ws4.Cells.Item(42, 6).Value = 4620
$ws4.Cells.Item(44, 6).Value = 3
$ws4.Cells.Item(47, 6).Value = 638
